# "Generate Report for Handback"
#
# The handback report workbook (Overview / zh-cn / de-de) is refreshed:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (Overview E/F columns, and the "Status" column on the zh-cn / de-de sheets)
#   - The "Latest Handback DateTime" for both locales is bumped to the new
#     handback timestamps
#   - The stale "Error Detail" note (out-of-date handback file warning) is
#     cleared now that the handback is in sync
#   - A couple of columns are widened on Overview / zh-cn / de-de so the new,
#     longer status text and file names are not clipped

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column (C), Latest Handback DateTime (K), Error Detail (P) ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K2").Value = "2016-11-14 07:30:16"
$zhcn.Range("K3").Value = "2016-11-14 07:30:16"
$zhcn.Range("P2").Value = ""

# --- de-de sheet: Status column (C), Latest Handback DateTime (K), Error Detail (P) ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("K2").Value = "2016-11-14 07:30:34"
$dede.Range("K3").Value = "2016-11-14 07:30:34"
$dede.Range("P2").Value = ""

# --- Column width adjustments (widen to fit the longer status text / error column) ---
# ColumnWidth is quantized internally to 1/6 character; these are the closest
# settings to the refreshed widths.
$overview.Columns.Item(5).ColumnWidth  = 29.166666666666668  # E
$overview.Columns.Item(6).ColumnWidth  = 29.166666666666668  # F

$zhcn.Columns.Item(3).ColumnWidth      = 29.166666666666668  # C (Status)
$zhcn.Columns.Item(16).ColumnWidth     = 12.833333333333334  # P (Error Detail)

$dede.Columns.Item(3).ColumnWidth      = 29.166666666666668  # C (Status)
$dede.Columns.Item(16).ColumnWidth     = 12.833333333333334  # P (Error Detail)
